$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.0463764283011649
$ws.Range("C2").Value = 0.1794809538289513

$ws.Range("B3").Value = 0.02157205171300387
$ws.Range("C3").Value = 0.09570174413718284

$ws.Range("B4").Value = 0.7519808397960515
$ws.Range("C4").Value = 0.3503770573730713

$ws.Range("B5").Value = 0.9169947986045117
$ws.Range("C5").Value = 0.7194900691260924

$ws.Range("B6").Value = 0.9792614230038085
$ws.Range("C6").Value = 0.8444421974125974

$ws.Range("B7").Value = 0.8281868950192681
$ws.Range("C7").Value = 0.31106226319429

$ws.Range("B8").Value = 0.01639171600341797
$ws.Range("C8").Value = 0.1298224544525146

$ws.Range("B9").Value = 0.01318404979094875
$ws.Range("C9").Value = 0.01383336160571978

$ws.Range("B10").Value = 0.01410578161886031
$ws.Range("C10").Value = 0.01830880674363285
